$wb = $excel.ActiveWorkbook
Write-Output "test"
